$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify steel description: drop the "RME/" token from the second line
$ws.Range("B2").Value = "36% CR+PC/LFM+CDL/H:1`n20% S/LFM+CDL/H:1`n13% S+SL/LFM+CDL/H:1`n6% S/LFBR+CDL/H:1`n15% CR/LFM+CDL/H:2`n2% CR/LFM+CDL/HBET:3-5`n8% MUR/LWAL+CDN/H:1"

# Wrap the (still multi-line) text and size row 2 to fit it
$ws.Range("B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 380

# Leave the selection as the author last left it
$ws.Range("B2:B12").Select() | Out-Null
